$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.970.39'
$ws.Range('E2').Value = '  +5.87%  '
$ws.Range('D3').Value = '2.746.22'
$ws.Range('E3').Value = '  +5.10%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.07'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +1.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.47'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +7.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.612'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +2.67%  '
$ws.Range('D9').Value = '2.783.49'
$ws.Range('E9').Value = '  +6.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.79'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +4.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.115'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +8.85%  '
$ws.Range('E12').Value = '  +4.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.159'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +1.77%  '
$ws.Range('D14').Value = '3.250.26'
$ws.Range('E14').Value = '  +5.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.75'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +8.55%  '
$ws.Range('D16').Value = '63.898.34'
$ws.Range('E16').Value = '  +5.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000154'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  +9.32%  '
$ws.Range('D18').Value = '2.771.95'
$ws.Range('E18').Value = '  +5.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.12'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +6.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.91'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +6.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '367.02'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +5.74%  '
$ws.Range('E22').Value = '  +2.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.542'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +1.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.19'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +4.18%  '
$ws.Range('E26').Value = '  +5.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.70'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').Value = '0.0₃0886'
$ws.Range('E29').Value = '  +11.10%  '
$ws.Range('E30').Value = '  +7.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.16'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +10.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '170.97'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.21'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +18.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '20.73'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +6.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.81'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +12.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.45'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +11.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.81'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +10.21%  '
$ws.Range('E39').Value = '  +20.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '354.26'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +11.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.28'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +9.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.27'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +2.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.71'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +13.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.46'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +12.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '146.16'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +7.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.29'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +11.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0597'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +8.51%  '
$ws.Range('E48').Value = '  +5.73%  '
$ws.Range('E49').Value = '  +7.37%  '
$ws.Range('E50').Value = '  +2.57%  '
$ws.Range('D51').Value = '2.180.05'
$ws.Range('E51').Value = '  +7.61%  '
